$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.4
$ws.Range("I2").Value = 2.75
$ws.Range("Y2").Value = 11
$ws.Range("AG2").Value = 151
$ws.Range("AH2").Value = 13
$ws.Range("AL2").Value = 21
$ws.Range("AO2").Value = 13
$ws.Range("AW2").Value = 401
$ws.Range("G3").Value = 2.6
$ws.Range("I3").Value = 2.5
$ws.Range("U3").Value = 1.4
$ws.Range("V3").Value = 2.75
$ws.Range("Z3").Value = 29
$ws.Range("AA3").Value = 19
$ws.Range("AD3").Value = 8.5
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 10
$ws.Range("Q6").Value = 2.08
$ws.Range("R6").Value = 1.73
$ws.Range("G7").Value = 1.44
$ws.Range("I7").Value = 7.5
$ws.Range("M7").Value = 1.07
$ws.Range("N7").Value = 8.5
$ws.Range("Q7").Value = 2.08
$ws.Range("R7").Value = 1.73
$ws.Range("X7").Value = 5.5
$ws.Range("Y7").Value = 9.5
$ws.Range("Z7").Value = 8.5
$ws.Range("AB7").Value = 41
$ws.Range("AD7").Value = 8.5
$ws.Range("AE7").Value = 26
$ws.Range("AU7").Value = 11
$ws.Range("G8").Value = 1.36
$ws.Range("H8").Value = 4.5
$ws.Range("I8").Value = 8.5
$ws.Range("J8").Value = 1.95
$ws.Range("L8").Value = 9
$ws.Range("R8").Value = 1.67
$ws.Range("U8").Value = 2.63
$ws.Range("V8").Value = 1.44
$ws.Range("Z8").Value = 8
$ws.Range("AC8").Value = 8.5
$ws.Range("AO8").Value = 6.5
$ws.Range("AX8").Value = 9.5
$ws.Range("AY8").Value = 51
$ws.Range("G9").Value = 1.67
$ws.Range("R9").Value = 1.6
$ws.Range("AA9").Value = 15
$ws.Range("AQ9").Value = 29
$ws.Range("AV9").Value = 67
$ws.Range("M10").Value = 1.1
$ws.Range("O10").Value = 1.5
$ws.Range("G11").Value = 2.3
$ws.Range("I11").Value = 2.8
$ws.Range("J11").Value = 3
$ws.Range("K11").Value = 2.25
$ws.Range("L11").Value = 3.4
$ws.Range("M11").Value = 1.04
$ws.Range("O11").Value = 1.22
$ws.Range("U11").Value = 1.62
$ws.Range("V11").Value = 2.2
$ws.Range("W11").Value = 9.5
$ws.Range("Y11").Value = 9.5
$ws.Range("Z11").Value = 23
$ws.Range("AC11").Value = 13
$ws.Range("AE11").Value = 12
$ws.Range("AI11").Value = 15
$ws.Range("AK11").Value = 29
$ws.Range("AL11").Value = 21
$ws.Range("AM11").Value = 26
$ws.Range("AO11").Value = 13
$ws.Range("AY11").Value = 15
$ws.Range("M12").Value = 1.04
$ws.Range("N12").Value = 12
$ws.Range("O12").Value = 1.22
$ws.Range("Q12").Value = 1.8
$ws.Range("R12").Value = 2
$ws.Range("G13").Value = 1.33
$ws.Range("H13").Value = 4.75
$ws.Range("J13").Value = 1.83
$ws.Range("K13").Value = 2.5
$ws.Range("L13").Value = 7.5
$ws.Range("M13").Value = 1.04
$ws.Range("N13").Value = 13
$ws.Range("U13").Value = 2
$ws.Range("V13").Value = 1.75
$ws.Range("Z13").Value = 8.5
$ws.Range("AB13").Value = 29
$ws.Range("AD13").Value = 9.5
$ws.Range("AF13").Value = 67
$ws.Range("AG13").Value = 401
$ws.Range("AK13").Value = 101
$ws.Range("AN13").Value = 3.25
$ws.Range("AU13").Value = 9.5
$ws.Range("AX13").Value = 9
$ws.Range("BC13").Value = 351
$ws.Range("Q14").Value = 1.67
$ws.Range("R14").Value = 2.15
$ws.Range("U14").Value = 1.73
$ws.Range("G15").Value = 2.45
$ws.Range("I15").Value = 2.45
$ws.Range("K15").Value = 2.3
$ws.Range("L15").Value = 3.1
$ws.Range("Q15").Value = 1.67
$ws.Range("R15").Value = 2.15
$ws.Range("X15").Value = 15
$ws.Range("AY15").Value = 13
$ws.Range("G16").Value = 1.85
$ws.Range("H16").Value = 3.1
$ws.Range("I16").Value = 4.1
$ws.Range("J16").Value = 2.63
$ws.Range("K16").Value = 2.05
$ws.Range("L16").Value = 4.75
$ws.Range("M16").Value = 1.08
$ws.Range("N16").Value = 8
$ws.Range("Q16").Value = 2.25
$ws.Range("R16").Value = 1.62
$ws.Range("S16").Value = 1.5
$ws.Range("T16").Value = 2.5
$ws.Range("W16").Value = 6.5
$ws.Range("X16").Value = 8.5
$ws.Range("AB16").Value = 34
$ws.Range("AC16").Value = 8
$ws.Range("AH16").Value = 10
$ws.Range("AK16").Value = 41
$ws.Range("AO16").Value = 11
$ws.Range("AQ16").Value = 41
$ws.Range("AR16").Value = 67
$ws.Range("AT16").Value = 2.5
$ws.Range("AY16").Value = 23
$ws.Range("U17").Value = 2.2
$ws.Range("V17").Value = 1.62
$ws.Range("X17").Value = 6.5
$ws.Range("AB17").Value = 29
$ws.Range("AN17").Value = 3.2
$ws.Range("AR17").Value = 34
$ws.Range("BA17").Value = 301
$ws.Range("BB17").Value = 251
$ws.Range("G18").Value = 3
$ws.Range("H18").Value = 3.2
$ws.Range("I18").Value = 2.4
$ws.Range("K18").Value = 2.1
$ws.Range("M18").Value = 1.06
$ws.Range("N18").Value = 9.5
$ws.Range("O18").Value = 1.3
$ws.Range("P18").Value = 3.4
$ws.Range("Q18").Value = 2.05
$ws.Range("R18").Value = 1.75
$ws.Range("U18").Value = 1.8
$ws.Range("V18").Value = 1.95
$ws.Range("AA18").Value = 23
$ws.Range("AC18").Value = 9.5
$ws.Range("AE18").Value = 13
$ws.Range("AH18").Value = 8
$ws.Range("AI18").Value = 12
$ws.Range("BC18").Value = 151
$ws.Range("G20").Value = 2.8
$ws.Range("H20").Value = 3.25
$ws.Range("J20").Value = 3.35
$ws.Range("K20").Value = 2.12
$ws.Range("L20").Value = 2.92
$ws.Range("O20").Value = 1.32
$ws.Range("P20").Value = 3.1
$ws.Range("Q20").Value = 1.95
$ws.Range("S20").Value = 1.38
$ws.Range("U20").Value = 1.75
$ws.Range("V20").Value = 1.98
$ws.Range("X20").Value = 14
$ws.Range("AH20").Value = 8
$ws.Range("AI20").Value = 12
$ws.Range("AL20").Value = 19.5
$ws.Range("AM20").Value = 29
$ws.Range("AN20").Value = 4.8
$ws.Range("AO20").Value = 15
$ws.Range("AQ20").Value = 70
$ws.Range("AU20").Value = 6.9
$ws.Range("AZ20").Value = 19
$ws.Range("BB20").Value = 75
$ws.Range("K21").Value = 2.5
$ws.Range("U21").Value = 2
$ws.Range("V21").Value = 1.75
$ws.Range("AC21").Value = 13
$ws.Range("AG21").Value = 401
$ws.Range("AJ21").Value = 9
$ws.Range("BA21").Value = 17
$ws.Range("BC21").Value = 126
$ws.Range("M22").Value = 1.05
$ws.Range("N22").Value = 8.5
$ws.Range("Q22").Value = 1.93
$ws.Range("R22").Value = 1.88
$ws.Range("U22").Value = 1.8
$ws.Range("V22").Value = 1.91
$ws.Range("M23").Value = 1.05
$ws.Range("N23").Value = 8.5
$ws.Range("U23").Value = 1.67
$ws.Range("M24").Value = 1.03
$ws.Range("N24").Value = 10
$ws.Range("Q24").Value = 1.67
$ws.Range("R24").Value = 2.15
$ws.Range("U24").Value = 1.91
$ws.Range("V24").Value = 1.8
$ws.Range("S25").Value = 1.28
$ws.Range("T25").Value = 3.5
